$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Copy() | Out-Null
$ws.Range("U2:W2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$ws.Range("U2:W2").Interior.ColorIndex = -4142

$ws.Range("U2").Value = "Method"
$ws.Range("V2").Value = "Massfrac"
$ws.Range("W2").Value = "Partdiam"

$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 0.04
$ws.Range("W3").Value = 0.0006
